$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remove highlight style from D1:E1 (back to default/Normal)
$ws1.Range("D1:E1").Style = "Normal"

# Add new column L
$ws1.Range("L1").Value = "alhic2302-155_2-stick2"
$ws1.Range("L2").Value = 38.9322
$ws1.Range("L3").Value = 39.5556
$ws1.Range("L4").Value = 39.8722
$ws1.Range("L5").Value = 37.795
$ws1.Range("L6").Value = 36.755
$ws1.Range("L7").Value = 37.586
$ws1.Range("L8").Value = 37.583
$ws1.Range("L9").Value = 38.054
$ws1.Range("L10").Value = 37.225
$ws1.Range("L11").Value = 35.913
$ws1.Range("L12").Value = 36.649
$ws1.Range("L13").Value = 35.966
$ws1.Range("L14").Value = 35.644
$ws1.Range("L15").Value = 33.734
$ws1.Range("L16").Value = 35.243
$ws1.Range("L17").Value = 31.893
$ws1.Range("L18").Select() | Out-Null

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A12").Value = "alhic2302-155_2-stick2"
$ws2.Range("B12").Value = 588.4
$ws2.Range("C12").Value = 126.4
$ws2.Range("D12").Formula = "=0+2+16"
$ws2.Range("E12").Value = 32
$ws2.Range("F12").Value = 0
$ws2.Range("G12").Value = 61
$ws2.Range("H12").Value = 31
$ws2.Range("I12").Value = "MUST ADD 16mm to offset to acount for lost section at top of stick"
$ws2.Range("J14").Select() | Out-Null
